# This workbook tracks flight arrivals at KRK. This edit appends 19 new
# arrival records (rows 462-480) pulled in via the "download from internet"
# panel, covering the tail end of Friday Jan 13 through Saturday Jan 14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one row: NUMBER, DATE, TIME, FLIGHT, FROM, SHORT,
# AIRLINE, MODEL, AIRCFAT ID, STATUS, DIFFERENCE (K and M stay blank).
$newRows = @(
    ,@(461.0, "Friday, Jan 13", "11:15 PM", "FR2362", "London", "(STN)", "Ryanair ", "B38M", "(SP-RZK)", "11:34 PM", $null, "0 hours, 19 minutes", $null)
    ,@(462.0, "Friday, Jan 13", "11:20 PM", "OS599", "Vienna", "(VIE)", "Austrian Airlines ", "E195", "(OE-LWP)", "11:09 PM", $null, "0 hours, -11 minutes", $null)
    ,@(463.0, "Friday, Jan 13", "11:25 PM", "LH1370", "Frankfurt", "(FRA)", "Lufthansa ", "A21N", "(D-AIEK)", "11:23 PM", $null, "0 hours, -2 minutes", $null)
    ,@(464.0, "Friday, Jan 13", "11:25 PM", "W65004", "London", "(LTN)", "Wizz Air ", "A21N", "(HA-LVG)", "11:13 PM", $null, "0 hours, -12 minutes", $null)
    ,@(465.0, "Friday, Jan 13", "11:30 PM", "FR3364", "Berlin", "(BER)", "Buzz ", "B38M", "(SP-RZF)", "11:19 PM", $null, "0 hours, -11 minutes", $null)
    ,@(466.0, "Friday, Jan 13", "11:30 PM", "W65042", "Bergen", "(BGO)", "Wizz Air ", "A321", "(HA-LXO)", "11:27 PM", $null, "0 hours, -3 minutes", $null)
    ,@(467.0, "Friday, Jan 13", "11:35 PM", "LO3911", "Warsaw", "(WAW)", "LOT ", "E190", "(SP-LMH)", "11:44 PM", $null, "0 hours, 9 minutes", $null)
    ,@(468.0, "Friday, Jan 13", "11:40 PM", "FR8673", "Girona", "(GRO)", "Ryanair ", "B738", "(SP-RSA)", "12:30 AM", $null, "0 hours, 50 minutes", $null)
    ,@(469.0, "Friday, Jan 13", "11:40 PM", "W65078", "Stockholm", "(NYO)", "Wizz Air ", "A21N", "(HA-LVO)", "11:32 PM", $null, "0 hours, -8 minutes", $null)
    ,@(470.0, "Friday, Jan 13", "11:59 PM", "FR6227", "Gran Canaria", "(LPA)", "Buzz ", "B38M", "(SP-RZH)", "12:26 AM", $null, "0 hours, 27 minutes", $null)
    ,@(471.0, "Saturday, Jan 14", "12:25 AM", "W65052", "Larnaca", "(LCA)", "Wizz Air ", "A21N", "(HA-LZI)", "11:51 PM", $null, "23 hours, 26 minutes", $null)
    ,@(472.0, "Saturday, Jan 14", "5:58 AM", "UNKNOWN", "Sofia", "(SOF)", "Ryanair ", "B738", "(SP-RKC)", "5:20 AM", $null, "0 hours, -38 minutes", $null)
    ,@(473.0, "Saturday, Jan 14", "6:24 AM", "UNKNOWN", "Katowice", "(KTW)", "Ryanair ", "B738", "(SP-RKB)", "6:30 AM", $null, "0 hours, 6 minutes", $null)
    ,@(474.0, "Saturday, Jan 14", "8:00 AM", "FR9662", "Rome", "(CIA)", "Ryanair ", "B738", "(9H-QAS)", "7:38 AM", $null, "0 hours, -22 minutes", $null)
    ,@(475.0, "Saturday, Jan 14", "8:05 AM", "AY1161", "Helsinki", "(HEL)", "Finnair ", "E190", "(OH-LKH)", "7:47 AM", $null, "0 hours, -18 minutes", $null)
    ,@(476.0, "Saturday, Jan 14", "8:36 AM", "E47904", "Antalya", "(AYT)", "Enter Air ", "B738", "(SP-ESH)", "8:13 AM", $null, "0 hours, -23 minutes", $null)
    ,@(477.0, "Saturday, Jan 14", "8:40 AM", "FR6319", "Marseille", "(MRS)", "Ryanair ", "B738", "(9H-QAC)", "8:20 AM", $null, "0 hours, -20 minutes", $null)
    ,@(478.0, "Saturday, Jan 14", "9:35 AM", "FR5623", "Edinburgh", "(EDI)", "Ryanair ", "B38M", "(EI-HGG)", "9:33 AM", $null, "0 hours, -2 minutes", $null)
    ,@(479.0, "Saturday, Jan 14", "9:45 AM", "LH1364", "Frankfurt", "(FRA)", "Lufthansa ", "A320", "(D-AIUK)", "9:45 AM", $null, "0 hours, 0 minutes", $null)
)

$startRow = 462
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($col = 1; $col -le $row.Count; $col++) {
        $val = $row[$col - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $col).Value() = $val
        }
    }
}
